$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.451.14"
$ws.Range("E2").Value = "  -0.12%  "

$ws.Range("D3").Value = "1.916.86"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("E4").Value = "  +0.63%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.45"
$ws.Range("E5").Value = "  -0.03%  "

$ws.Range("E6").Value = "  +0.70%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4809"
$ws.Range("E7").Value = "  -0.69%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4048"
$ws.Range("E8").Value = "  -0.34%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08205"
$ws.Range("E9").Value = "  +1.03%  "

$ws.Range("E10").Value = "  -0.10%  "

$ws.Range("E11").Value = "  -0.68%  "

$ws.Range("D12").Value = "1.903.34"
$ws.Range("E12").Value = "  -1.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.048"
$ws.Range("E13").Value = "  +0.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.221"
$ws.Range("E14").Value = "  +1.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.16"
$ws.Range("E15").Value = "  +0.83%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06852"
$ws.Range("E16").Value = "  +1.74%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.012"
$ws.Range("E17").Value = "  +0.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001039"
$ws.Range("E18").Value = "  -0.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.48"
$ws.Range("E19").Value = "  -1.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.011"
$ws.Range("E20").Value = "  +0.59%  "

$ws.Range("D21").Value = "29.465.99"
$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.663"
$ws.Range("E22").Value = "  +1.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.83"
$ws.Range("E23").Value = "  +0.38%  "

$ws.Range("E24").Value = "  +1.51%  "

$ws.Range("D25").Value = "2.139.74"
$ws.Range("E25").Value = "  -0.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.609"
$ws.Range("E26").Value = "  +5.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.84"
$ws.Range("E27").Value = "  +1.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.00"
$ws.Range("E28").Value = "  -0.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.095"
$ws.Range("E29").Value = "  -0.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.38"
$ws.Range("E30").Value = "  +1.24%  "

$ws.Range("E31").Value = "  -2.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09609"
$ws.Range("E32").Value = "  +0.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.603"
$ws.Range("E33").Value = "  +1.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.560"
$ws.Range("E34").Value = "  +0.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.369"
$ws.Range("E35").Value = "  -1.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06289"
$ws.Range("E36").Value = "  +3.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02284"
$ws.Range("E37").Value = "  +0.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.181"
$ws.Range("E38").Value = "  +0.68%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5923"
$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.71"
$ws.Range("E40").Value = "  +3.59%  "

$ws.Range("E41").Value = "  +0.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.863"
$ws.Range("E42").Value = "  -0.65%  "

$ws.Range("E43").Value = "  -0.58%  "

$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.391"
$ws.Range("E45").Value = "  -0.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.40"
$ws.Range("E46").Value = "  +0.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07469"
$ws.Range("E47").Value = "  -3.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5555"
$ws.Range("E48").Value = "  +0.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.928"
$ws.Range("E49").Value = "  -1.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "117.96"
$ws.Range("E50").Value = "  +2.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.422"
$ws.Range("E51").Value = "  +3.29%  "

